# repull data, push all data, mean calculation
# Update column F (dSF) values to reflect the re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 4
